$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "34.521.32"
$ws.Range("E2").Value = "  +1.17%  "
$ws.Range("D3").Value = "1.787.57"
$ws.Range("E3").Value = "  -0.18%  "
$ws.Range("E4").Value = "  +0.22%  "
$ws.Range("D5").Value = "222.63"
$ws.Range("E5").Value = "  -1.74%  "
$ws.Range("E6").Value = "  -1.28%  "
$ws.Range("E7").Value = "  +0.20%  "
$ws.Range("E8").Value = "  +5.93%  "
$ws.Range("D9").Value = "0.280"
$ws.Range("E9").Value = "  +0.03%  "
$ws.Range("D10").Value = "0.0687"
$ws.Range("E10").Value = "  +2.79%  "
$ws.Range("E11").Value = "  +1.24%  "
$ws.Range("D12").Value = "2.044.90"
$ws.Range("E12").Value = "  -0.15%  "
$ws.Range("D13").Value = "1.788.53"
$ws.Range("E13").Value = "  -0.35%  "
$ws.Range("D14").Value = "10.84"
$ws.Range("E14").Value = "  +5.19%  "
$ws.Range("D15").Value = "34.508.86"
$ws.Range("E15").Value = "  +1.10%  "
$ws.Range("D16").Value = "0.630"
$ws.Range("E16").Value = "  +0.37%  "
$ws.Range("D17").Value = "4.27"
$ws.Range("E17").Value = "  +1.67%  "
$ws.Range("D18").Value = "68.62"
$ws.Range("E18").Value = "  -0.50%  "
$ws.Range("D19").Value = "253.76"
$ws.Range("E19").Value = "  +0.49%  "
$ws.Range("D20").Value = "0.0₃0782"
$ws.Range("E20").Value = "  +5.36%  "
$ws.Range("E21").Value = "  +0.28%  "
$ws.Range("D22").Value = "10.43"
$ws.Range("E22").Value = "  +0.66%  "
$ws.Range("D23").Value = "4.15"
$ws.Range("E23").Value = "  -1.80%  "
$ws.Range("E24").Value = "  +0.39%  "
$ws.Range("D25").Value = "160.18"
$ws.Range("E25").Value = "  +1.03%  "
$ws.Range("D26").Value = "16.37"
$ws.Range("E26").Value = "  -0.96%  "
$ws.Range("D27").Value = "7.07"
$ws.Range("E27").Value = "  +0.82%  "
$ws.Range("E28").Value = "  -1.24%  "
$ws.Range("E29").Value = "  +0.37%  "
$ws.Range("D30").Value = "0.0515"
$ws.Range("E30").Value = "  +0.04%  "
$ws.Range("D31").Value = "3.75"
$ws.Range("E31").Value = "  -2.37%  "
$ws.Range("E32").Value = "  -0.76%  "
$ws.Range("D33").Value = "3.55"
$ws.Range("E33").Value = "  -1.14%  "
$ws.Range("D34").Value = "1.87"
$ws.Range("E34").Value = "  +0.57%  "
$ws.Range("D35").Value = "1.434.35"
$ws.Range("E35").Value = "  -4.67%  "
$ws.Range("D36").Value = "0.636"
$ws.Range("E36").Value = "  -0.16%  "
$ws.Range("E37").Value = "  -1.55%  "
$ws.Range("E38").Value = "  +2.19%  "
$ws.Range("D39").Value = "84.22"
$ws.Range("E39").Value = "  +0.70%  "
$ws.Range("D40").Value = "2.84"
$ws.Range("E40").Value = "  +4.49%  "
$ws.Range("D41").Value = "2.36"
$ws.Range("E41").Value = "  +0.48%  "
$ws.Range("D42").Value = "0.915"
$ws.Range("E42").Value = "  +1.08%  "
$ws.Range("D43").Value = "2.07"
$ws.Range("E43").Value = "  +1.19%  "
$ws.Range("D44").Value = "5.96"
$ws.Range("E44").Value = "  +3.84%  "
$ws.Range("E45").Value = "  -1.15%  "
$ws.Range("E46").Value = "  -4.38%  "
$ws.Range("D47").Value = "1.946.04"
$ws.Range("E47").Value = "  +0.03%  "
$ws.Range("B48").Value = "PaxDollar"
$ws.Range("C48").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D48").Value = "1.00"
$ws.Range("E48").Value = "  +0.25%  "
$ws.Range("D49").Value = "11.96"
$ws.Range("E49").Value = "  +1.16%  "
$ws.Range("B50").Value = "Quant"
$ws.Range("C50").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D50").Value = "103.67"
$ws.Range("E50").Value = "  +5.38%  "
$ws.Range("D51").Value = "49.69"
$ws.Range("E51").Value = "  -3.54%  "
